# LOM3203.xlsx edit: re-shuffle/correct the "Programa"/"Avaliacao" block
# content (rows 10, 13-25) per the source commit, and drop the now-unused
# trailing row. Existing text is moved with Range.Copy so the destination
# cell inherits the correct shared-string/number-format/style behaviour
# (avoids Excel's automatic date-detection on "01/01/2012"-like strings).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- move reused text first, before the source cells are overwritten ----

# row10 B/C <- row13 B/C ("5982760 - Carlos Alberto Baldan")
$ws.Range("B13").Copy($ws.Range("B10"))
$ws.Range("C13").Copy($ws.Range("C10"))

# row15 B/C <- row8 B/C ("01/01/2012")
$ws.Range("B8").Copy($ws.Range("B15"))
$ws.Range("C8").Copy($ws.Range("C15"))

# row18 B/C <- row13 B/C ("5982760 - Carlos Alberto Baldan", still intact)
$ws.Range("B13").Copy($ws.Range("B18"))
$ws.Range("C13").Copy($ws.Range("C18"))

# row23 B/C <- row24 B/C ("LOM3206 -  Eletrônica  (Requisito)")
$ws.Range("B24").Copy($ws.Range("B23"))
$ws.Range("C24").Copy($ws.Range("C23"))

# row24 B/C <- row25 B/C ("LOM3260 -  Computação Científica em Python  (Requisito)")
$ws.Range("B25").Copy($ws.Range("B24"))
$ws.Range("C25").Copy($ws.Range("C24"))

# ---- now the label (column A) and the few genuinely-new values ----------

# row13: "Programa resumido:" / "Semestral"
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("A13").Font.Bold = $true
$ws.Range("A13").WrapText = $false
$ws.Range("A13").VerticalAlignment = -4160
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# row14: just the "Short syllabus:" label now
$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("B14").ClearContents()
$ws.Range("C14").ClearContents()

# row15: "Programa:" label (data already copied above)
$ws.Range("A15").Value = "Programa:"

# row16: just the "Syllabus:" label now
$ws.Range("A16").Value = "Syllabus:"
$ws.Range("B16").ClearContents()
$ws.Range("C16").ClearContents()

# row17: "Avaliação:" label
$ws.Range("A17").Value = "Avaliação:"

# row18: "Método:" label (data already copied above)
$ws.Range("A18").Value = "Método:"

# row19: "Critério:" label (B/C untouched)
$ws.Range("A19").Value = "Critério:"

# row20: "Norma de recuperação:" label (B/C untouched)
$ws.Range("A20").Value = "Norma de recuperação:"

# row21: "Bibliografia:" label (B/C untouched)
$ws.Range("A21").Value = "Bibliografia:"

# row22: just the "Requisitos:" label now
$ws.Range("A22").Value = "Requisitos:"
$ws.Range("B22").ClearContents()
$ws.Range("C22").ClearContents()

# row23: no label any more (data already copied above)
$ws.Range("A23").ClearContents()

# ---- row 25 is now unused; drop it entirely ------------------------------
$ws.Rows(25).Delete()

# ---- row heights ----------------------------------------------------------
$ws.Rows(13).RowHeight = 60
$ws.Rows(14).RowHeight = 60
$ws.Rows(15).RowHeight = 120
$ws.Rows(16).RowHeight = 120
$ws.Rows(17).AutoFit()
$ws.Rows(18).RowHeight = 60
$ws.Rows(21).RowHeight = 120
$ws.Rows(22).AutoFit()
$ws.Rows(23).RowHeight = 30
